# Edit script for cs-en-us-pbbs.xlsx weekly crime data update.
# - Updates the report volume/issue number and the covered date range
#   (rich-text header cells A8 and C9), using Characters() so we only
#   touch the digits that changed (right-to-left so earlier offsets
#   stay valid after the date length changes from 8 to 9 chars).
# - Updates all the weekly/28-day/YTD/2-year crime statistic cells in
#   the precinct table (rows 14-30, columns C-N) to the newly reported
#   figures and their recomputed percentage changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Volume 30   Number  23" -> "...24"
$ws.Range("A8").Characters(21, 2).Text = "24"

# --- Header: "Report Covering the Week  6/5/2023  Through  6/11/2023"
#     -> "...6/12/2023  Through  6/18/2023"
# (Second date replaced first so the first date's offset isn't shifted
#  by the length change of "6/5/2023" (8 chars) -> "6/12/2023" (9 chars).)
$ws.Range("C9").Characters(46, 9).Text = "6/18/2023"
$ws.Range("C9").Characters(27, 8).Text = "6/12/2023"

# --- Precinct weekly crime statistics table (rows 14-30)
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 5
$ws.Range("H14").Value = -60
$ws.Range("I14").Value = 31
$ws.Range("J14").Value = 29
$ws.Range("K14").Value = 6.896551724137
$ws.Range("L14").Value = 19.230769230769
$ws.Range("M14").Value = -8.823529411764
$ws.Range("N14").Value = -71.028037383177
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 7
$ws.Range("E15").Value = -42.857142857142
$ws.Range("F15").Value = 16
$ws.Range("G15").Value = 20
$ws.Range("H15").Value = -20
$ws.Range("I15").Value = 97
$ws.Range("J15").Value = 103
$ws.Range("K15").Value = -5.825242718446
$ws.Range("L15").Value = 14.117647058823
$ws.Range("M15").Value = 36.619718309859
$ws.Range("N15").Value = -61.507936507936
$ws.Range("C16").Value = 30
$ws.Range("D16").Value = 38
$ws.Range("E16").Value = -21.052631578947
$ws.Range("F16").Value = 132
$ws.Range("G16").Value = 148
$ws.Range("H16").Value = -10.810810810810
$ws.Range("I16").Value = 800
$ws.Range("J16").Value = 877
$ws.Range("K16").Value = -8.779931584948
$ws.Range("L16").Value = 39.860139860139
$ws.Range("M16").Value = -35.431799838579
$ws.Range("N16").Value = -87.505856629704
$ws.Range("C17").Value = 75
$ws.Range("D17").Value = 81
$ws.Range("E17").Value = -7.407407407407
$ws.Range("F17").Value = 294
$ws.Range("G17").Value = 284
$ws.Range("H17").Value = 3.521126760563
$ws.Range("I17").Value = 1571
$ws.Range("J17").Value = 1510
$ws.Range("K17").Value = 4.039735099337
$ws.Range("L17").Value = 23.215686274509
$ws.Range("M17").Value = 47.928436911487
$ws.Range("N17").Value = -47.458193979933
$ws.Range("C18").Value = 21
$ws.Range("D18").Value = 35
$ws.Range("E18").Value = -40
$ws.Range("F18").Value = 108
$ws.Range("G18").Value = 151
$ws.Range("H18").Value = -28.476821192053
$ws.Range("I18").Value = 812
$ws.Range("J18").Value = 900
$ws.Range("K18").Value = -9.777777777777
$ws.Range("L18").Value = 12.621359223301
$ws.Range("M18").Value = -45.794392523364
$ws.Range("N18").Value = -90.471720253461
$ws.Range("C19").Value = 126
$ws.Range("D19").Value = 135
$ws.Range("E19").Value = -6.666666666666
$ws.Range("F19").Value = 483
$ws.Range("G19").Value = 519
$ws.Range("H19").Value = -6.936416184971
$ws.Range("I19").Value = 3011
$ws.Range("J19").Value = 3201
$ws.Range("K19").Value = -5.935645110902
$ws.Range("L19").Value = 46.949731576378
$ws.Range("M19").Value = 26.939291736930
$ws.Range("N19").Value = -23.849266565503
$ws.Range("C20").Value = 41
$ws.Range("D20").Value = 34
$ws.Range("E20").Value = 20.588235294117
$ws.Range("F20").Value = 147
$ws.Range("G20").Value = 142
$ws.Range("H20").Value = 3.521126760563
$ws.Range("I20").Value = 812
$ws.Range("J20").Value = 794
$ws.Range("K20").Value = 2.267002518891
$ws.Range("L20").Value = 61.111111111111
$ws.Range("M20").Value = -10.077519379845
$ws.Range("N20").Value = -92.696528152545
$ws.Range("C21").Value = 298
$ws.Range("D21").Value = 331
$ws.Range("E21").Value = -9.969788519637
$ws.Range("F21").Value = 1182
$ws.Range("G21").Value = 1269
$ws.Range("H21").Value = -6.855791962174
$ws.Range("I21").Value = 7134
$ws.Range("J21").Value = 7414
$ws.Range("K21").Value = -3.776638791475
$ws.Range("L21").Value = 36.353211009174
$ws.Range("M21").Value = -0.626828249059
$ws.Range("N21").Value = -78.606129670725
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 300
$ws.Range("F22").Value = 16
$ws.Range("G22").Value = 17
$ws.Range("H22").Value = -5.882352941176
$ws.Range("I22").Value = 78
$ws.Range("J22").Value = 91
$ws.Range("K22").Value = -14.285714285714
$ws.Range("L22").Value = 21.875
$ws.Range("M22").Value = -37.6
$ws.Range("D23").Value = 20
$ws.Range("E23").Value = -60
$ws.Range("F23").Value = 43
$ws.Range("G23").Value = 59
$ws.Range("H23").Value = -27.118644067796
$ws.Range("I23").Value = 224
$ws.Range("J23").Value = 233
$ws.Range("K23").Value = -3.862660944206
$ws.Range("L23").Value = 3.225806451612
$ws.Range("M23").Value = 73.643410852713
$ws.Range("C24").Value = 332
$ws.Range("D24").Value = 318
$ws.Range("E24").Value = 4.402515723270
$ws.Range("F24").Value = 1233
$ws.Range("G24").Value = 1205
$ws.Range("H24").Value = 2.323651452282
$ws.Range("I24").Value = 7434
$ws.Range("J24").Value = 6913
$ws.Range("K24").Value = 7.536525386952
$ws.Range("L24").Value = 44.377549038648
$ws.Range("M24").Value = 38.875396973659
$ws.Range("C25").Value = 129
$ws.Range("E25").Value = -3.007518796992
$ws.Range("F25").Value = 488
$ws.Range("G25").Value = 517
$ws.Range("H25").Value = -5.609284332688
$ws.Range("I25").Value = 2647
$ws.Range("J25").Value = 2513
$ws.Range("K25").Value = 5.332272184639
$ws.Range("L25").Value = 24.917413874469
$ws.Range("M25").Value = -13.946684005201
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -20
$ws.Range("F26").Value = 23
$ws.Range("G26").Value = 29
$ws.Range("H26").Value = -20.689655172413
$ws.Range("I26").Value = 144
$ws.Range("J26").Value = 161
$ws.Range("K26").Value = -10.559006211180
$ws.Range("L26").Value = -1.369863013698
$ws.Range("C27").Value = 10
$ws.Range("D27").Value = 15
$ws.Range("E27").Value = -33.333333333333
$ws.Range("F27").Value = 53
$ws.Range("G27").Value = 62
$ws.Range("H27").Value = -14.516129032258
$ws.Range("I27").Value = 284
$ws.Range("J27").Value = 317
$ws.Range("K27").Value = -10.410094637224
$ws.Range("L27").Value = 14.056224899598
$ws.Range("C28").Value = 2
$ws.Range("E28").Value = -33.333333333333
$ws.Range("F28").Value = 10
$ws.Range("G28").Value = 13
$ws.Range("H28").Value = -23.076923076923
$ws.Range("I28").Value = 69
$ws.Range("J28").Value = 90
$ws.Range("K28").Value = -23.333333333333
$ws.Range("L28").Value = -25.806451612903
$ws.Range("M28").Value = -31
$ws.Range("N28").Value = -79.824561403508
$ws.Range("C29").Value = 2
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = -33.333333333333
$ws.Range("F29").Value = 8
$ws.Range("G29").Value = 11
$ws.Range("H29").Value = -27.272727272727
$ws.Range("I29").Value = 58
$ws.Range("J29").Value = 68
$ws.Range("K29").Value = -14.705882352941
$ws.Range("L29").Value = -30.120481927710
$ws.Range("M29").Value = -28.395061728395
$ws.Range("N29").Value = -80.983606557377
$ws.Range("F30").Value = 8
$ws.Range("G30").Value = 10
$ws.Range("H30").Value = -20
$ws.Range("I30").Value = 35
$ws.Range("J30").Value = 63
$ws.Range("K30").Value = -44.444444444444
$ws.Range("L30").Value = 9.375
